$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10
$ws.Range("C10").Style = "Bad"
$ws.Range("C10").Value = "NO"
$ws.Range("D10").Value = "Partially rebuilt, difficulties implementing unit tests"
$ws.Rows.Item(10).RowHeight = 15

# Row 11
$ws.Range("C11").Style = "Bad"
$ws.Range("C11").Value = "NO"
$ws.Range("D11").Value = "Only moving monitors implemented"
$ws.Rows.Item(11).RowHeight = 15

# Row 12
$ws.Range("C12").Style = "Bad"
$ws.Range("C12").Value = "NO"
$ws.Range("D12").Value = "Further discussion with Luxsonic required"
$ws.Rows.Item(12).RowHeight = 15

# Row 13
$ws.Range("C13").Style = "Bad"
$ws.Range("C13").Value = "NO"
$ws.Rows.Item(13).RowHeight = 15

# Row 14
$ws.Range("C14").Style = "Good"
$ws.Range("C14").Value = "YES"
$ws.Range("D14").Value = "Integrated by 2/12/17, presented at client meeting"
$ws.Rows.Item(14).RowHeight = 15

# Row 15
$ws.Range("A15").Value = 42783
$ws.Range("A8").Copy()
$ws.Range("A15").PasteSpecial(-4122)

# Update selection to D15
$ws.Range("D15").Select()
